$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New development log entry on row 12: copy the formatting from the row
# above (row 11) so the new row matches the existing table styling
# (date format, time format, borders, etc.), then fill in the values.
$ws.Range("C11:G11").Copy()
$ws.Range("C12:G12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C12").Value = 44630
$ws.Range("D12").Value = 0.52083333333333337
$ws.Range("E12").Value = 1.25
$ws.Range("F12").Value = "OO Design + Implementation"
$ws.Range("G12").Value = "OO Design + Implementation"

$ws.Range("E15").Select()
